$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1274
$ws.Range("F5").Value = 10
$ws.Range("F7").Value = 836
$ws.Range("F8").Value = 33
$ws.Range("F9").Value = 6713
$ws.Range("F11").Value = 95
$ws.Range("F12").Value = 134
$ws.Range("F13").Value = 6379
$ws.Range("F14").Value = 120
$ws.Range("F16").Value = 4299
$ws.Range("F19").Value = 4249
$ws.Range("F20").Value = 216
$ws.Range("F21").Value = 223
$ws.Range("F23").Value = 302
$ws.Range("F27").Value = 163
$ws.Range("F29").Value = 60
$ws.Range("F31").Value = 66
$ws.Range("F32").Value = 7747
$ws.Range("F34").Value = 1301
$ws.Range("F35").Value = 639
$ws.Range("F36").Value = 13
$ws.Range("F38").Value = 982
$ws.Range("F39").Value = 67
$ws.Range("F40").Value = 1541
$ws.Range("F42").Value = 882
$ws.Range("F43").Value = 37
$ws.Range("F44").Value = 3838
$ws.Range("F45").Value = 343
$ws.Range("F48").Value = 823
$ws.Range("F49").Value = 1070

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 10
$ws.Range("F16").Value = 74
$ws.Range("F19").Value = 864

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 231

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 231
$ws.Range("F7").Value = 1274
$ws.Range("F8").Value = 10
$ws.Range("F10").Value = 10
$ws.Range("F11").Value = 836
$ws.Range("F12").Value = 33
$ws.Range("F13").Value = 6713
$ws.Range("F15").Value = 95
$ws.Range("F16").Value = 134
$ws.Range("F17").Value = 6379
$ws.Range("F18").Value = 120
$ws.Range("F20").Value = 4299
$ws.Range("F21").Value = 4249
$ws.Range("F22").Value = 216
$ws.Range("F23").Value = 223
$ws.Range("F25").Value = 302
$ws.Range("F29").Value = 60
$ws.Range("F31").Value = 7747
$ws.Range("F33").Value = 1301
$ws.Range("F34").Value = 639
$ws.Range("F35").Value = 13
$ws.Range("F37").Value = 982
$ws.Range("F38").Value = 67
$ws.Range("F39").Value = 1541
$ws.Range("F41").Value = 882
$ws.Range("F42").Value = 37
$ws.Range("F43").Value = 3838
$ws.Range("F44").Value = 343
$ws.Range("F46").Value = 105
$ws.Range("F47").Value = 823
$ws.Range("F48").Value = 1070

